$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add two new rows of data (row 65 and 66)
$ws.Range("A65").Value = 44182
$ws.Range("B65").Value = 1440
$ws.Range("C65").Value = 379
$ws.Range("D65").Value = 1819

$ws.Range("A66").Value = 44183
$ws.Range("B66").Value = 1510
$ws.Range("C66").Value = 387
$ws.Range("D66").Value = 1897

# Update the view: scroll position and selection
$ws.Range("C47").Select()
$excel.ActiveWindow.ScrollRow = 38
$excel.ActiveWindow.ScrollColumn = 1
